$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G2: FechaSiniestro value changes from 29/04/2022 to 26/04/2021.
# Keep the leading apostrophe so the value stays text (quote-prefixed) like
# the original cell, instead of Excel auto-converting it to a date serial.
$ws.Range("G2").Value = "'26/04/2021"

# Update E2: NroPoliza value changes from numeric 12112002243 to text "12112002294".
# The leading apostrophe forces text-with-quote-prefix entry (same as typing
# '12112002294 into the cell), which is how the NroPoliza became a shared
# string instead of a plain number.
$ws.Range("E2").Value = "'12112002294"

# Update selection to D10
$ws.Range("D10").Select()
